# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) for the affected
# leve rows across sheets, matching the latest market snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 525
$ws.Range("I6").Value = 525
$ws.Range("K6").Value = 1575
$ws.Range("M6").Value = -1463

$ws.Range("H28").Value = 794303.4
$ws.Range("I28").Value = 1389185.1
$ws.Range("J28").Value = 1127.6666
$ws.Range("K28").Value = 1389185.1
$ws.Range("L28").Value = 1127.6666
$ws.Range("M28").Value = -1388700.1
$ws.Range("N28").Value = -2097.6666

$ws.Range("H58").Value = 77.5
$ws.Range("I58").Value = 77.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 232.5
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -82.5
$ws.Range("N58").ClearContents()

$ws.Range("H62").Value = 5924.7856
$ws.Range("I62").Value = 4263.222
$ws.Range("K62").Value = 4263.222
$ws.Range("M62").Value = -3639.222

$ws.Range("H65").Value = 5924.7856
$ws.Range("I65").Value = 4263.222
$ws.Range("K65").Value = 21316.11
$ws.Range("M65").Value = -18196.11

$ws.Range("H87").Value = 26478.572
$ws.Range("J87").Value = 26478.572
$ws.Range("L87").Value = 26478.572
$ws.Range("N87").Value = -28974.572

$ws.Range("H90").Value = 26478.572
$ws.Range("J90").Value = 26478.572
$ws.Range("L90").Value = 79435.716
$ws.Range("N90").Value = -91915.716

$ws.Range("H107").Value = 741501.7
$ws.Range("I107").Value = 1587799.4
$ws.Range("J107").Value = 991.125
$ws.Range("K107").Value = 1587799.4
$ws.Range("L107").Value = 991.125
$ws.Range("M107").Value = -1585879.4
$ws.Range("N107").Value = -4831.125

$ws.Range("H135").Value = 14577.875
$ws.Range("I135").Value = 14577.875
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 131200.875
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -128665.875
$ws.Range("N135").ClearContents()

$ws.Range("H137").Value = 31251568
$ws.Range("I137").Value = 40001428
$ws.Range("J137").Value = 2071.8572
$ws.Range("K137").Value = 120004284
$ws.Range("L137").Value = 6215.571599999999
$ws.Range("M137").Value = -120001734
$ws.Range("N137").Value = -11315.5716

$ws.Range("H138").Value = 6050678.5
$ws.Range("I138").Value = 1403441
$ws.Range("J138").Value = 7695086
$ws.Range("K138").Value = 4210323
$ws.Range("L138").Value = 23085258
$ws.Range("M138").Value = -4205183
$ws.Range("N138").Value = -23095538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8113.143
$ws.Range("I2").Value = 17500
$ws.Range("J2").Value = 1073
$ws.Range("K2").Value = 17500
$ws.Range("L2").Value = 1073
$ws.Range("M2").Value = -17387
$ws.Range("N2").Value = -1299

$ws.Range("H32").Value = 24540.791
$ws.Range("I32").Value = 2824.724
$ws.Range("J32").Value = 57686.367
$ws.Range("K32").Value = 2824.724
$ws.Range("L32").Value = 57686.367
$ws.Range("M32").Value = -2537.724
$ws.Range("N32").Value = -58260.367

$ws.Range("H110").Value = 2002.5
$ws.Range("I110").Value = 2002.5
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 2002.5
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 42.5
$ws.Range("N110").ClearContents()

$ws.Range("H116").Value = 8113.143
$ws.Range("I116").Value = 17500
$ws.Range("J116").Value = 1073
$ws.Range("K116").Value = 17500
$ws.Range("L116").Value = 1073
$ws.Range("M116").Value = -15206
$ws.Range("N116").Value = -5661

$ws.Range("H123").Value = 39617.332
$ws.Range("J123").Value = 39617.332
$ws.Range("L123").Value = 39617.332
$ws.Range("N123").Value = -49417.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8113.143
$ws.Range("I3").Value = 17500
$ws.Range("J3").Value = 1073
$ws.Range("K3").Value = 17500
$ws.Range("L3").Value = 1073
$ws.Range("M3").Value = -17386
$ws.Range("N3").Value = -1301

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

$ws.Range("H58").Value = 3116.8333
$ws.Range("I58").Value = 1694.4
$ws.Range("J58").Value = 4132.857
$ws.Range("K58").Value = 1694.4
$ws.Range("L58").Value = 4132.857
$ws.Range("M58").Value = -1491.4
$ws.Range("N58").Value = -4538.857

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

$ws.Range("H136").Value = 3116.8333
$ws.Range("I136").Value = 1694.4
$ws.Range("J136").Value = 4132.857
$ws.Range("K136").Value = 5083.200000000001
$ws.Range("L136").Value = 12398.571
$ws.Range("M136").Value = -2533.200000000001
$ws.Range("N136").Value = -17498.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 699.73334
$ws.Range("I121").Value = 170
$ws.Range("J121").Value = 832.1667
$ws.Range("K121").Value = 510
$ws.Range("L121").Value = 2496.5001
$ws.Range("M121").Value = 800
$ws.Range("N121").Value = -5116.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 484.625
$ws.Range("I107").Value = 393.33334
$ws.Range("J107").Value = 539.4
$ws.Range("K107").Value = 393.33334
$ws.Range("L107").Value = 539.4
$ws.Range("M107").Value = 1526.66666
$ws.Range("N107").Value = -4379.4

$ws.Range("H113").Value = 1083.3334
$ws.Range("I113").Value = 1083.3334
$ws.Range("K113").Value = 1083.3334
$ws.Range("M113").Value = 1086.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 506.92307
$ws.Range("I107").Value = 509.33334
$ws.Range("J107").Value = 501.5
$ws.Range("K107").Value = 1528.00002
$ws.Range("L107").Value = 1504.5
$ws.Range("M107").Value = 391.9999800000001
$ws.Range("N107").Value = -5344.5

$ws.Range("H123").Value = 34528.1
$ws.Range("J123").Value = 34528.1
$ws.Range("L123").Value = 34528.1
$ws.Range("N123").Value = -44328.1

$ws.Range("H132").Value = 11631075
$ws.Range("J132").Value = 2444.7058
$ws.Range("L132").Value = 7334.117400000001
$ws.Range("N132").Value = -12394.1174

$ws.Range("H136").Value = 18576280
$ws.Range("I136").Value = 27862236
$ws.Range("K136").Value = 83586708
$ws.Range("M136").Value = -83584158
